$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.388.08'
$ws.Cells.Item(2, 5).Value = '  -0.47%  '
$ws.Cells.Item(3, 4).Value = '1.845.59'
$ws.Cells.Item(3, 5).Value = '  -0.28%  '
$ws.Cells.Item(4, 4).Value = '0.9987'
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 4).Value = '240.64'
$ws.Cells.Item(5, 5).Value = '  -1.01%  '
$ws.Cells.Item(6, 4).Value = '0.6387'
$ws.Cells.Item(6, 5).Value = '  +0.24%  '
$ws.Cells.Item(7, 4).Value = '0.9999'
$ws.Cells.Item(7, 5).Value = '  +0.00%  '
$ws.Cells.Item(8, 2).Value = 'WrappedEther'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(8, 4).Value = '3.562.17'
$ws.Cells.Item(8, 5).Value = '  +91.62%  '
$ws.Cells.Item(9, 2).Value = 'Dogecoin'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(9, 4).Value = '0.07563'
$ws.Cells.Item(9, 5).Value = '  -0.11%  '
$ws.Cells.Item(10, 2).Value = 'Cardano'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Cells.Item(10, 4).Value = '0.2967'
$ws.Cells.Item(10, 5).Value = '  -0.95%  '
$ws.Cells.Item(11, 2).Value = 'Solana'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(11, 4).Value = '24.76'
$ws.Cells.Item(11, 5).Value = '  +1.90%  '
$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(12, 4).Value = '0.07742'
$ws.Cells.Item(12, 5).Value = '  +0.59%  '
$ws.Cells.Item(13, 4).Value = '4.994'
$ws.Cells.Item(13, 5).Value = '  -0.85%  '
$ws.Cells.Item(14, 4).Value = '0.6849'
$ws.Cells.Item(14, 5).Value = '  -0.45%  '
$ws.Cells.Item(15, 4).Value = '83.25'
$ws.Cells.Item(15, 5).Value = '  -1.01%  '
$ws.Cells.Item(16, 4).Value = '0.000009953'
$ws.Cells.Item(16, 5).Value = '  +1.59%  '
$ws.Cells.Item(17, 4).Value = '6.182'
$ws.Cells.Item(17, 5).Value = '  -1.54%  '
$ws.Cells.Item(18, 4).Value = '29.413.71'
$ws.Cells.Item(18, 5).Value = '  -0.49%  '
$ws.Cells.Item(19, 4).Value = '229.79'
$ws.Cells.Item(19, 5).Value = '  -3.27%  '
$ws.Cells.Item(20, 4).Value = '12.47'
$ws.Cells.Item(20, 5).Value = '  -0.55%  '
$ws.Cells.Item(21, 4).Value = '0.9999'
$ws.Cells.Item(21, 5).Value = '  -0.02%  '
$ws.Cells.Item(22, 4).Value = '7.576'
$ws.Cells.Item(22, 5).Value = '  -0.74%  '
$ws.Cells.Item(23, 4).Value = '1.000'
$ws.Cells.Item(23, 5).Value = '  +0.02%  '
$ws.Cells.Item(24, 4).Value = '156.52'
$ws.Cells.Item(24, 5).Value = '  +0.00%  '
$ws.Cells.Item(25, 4).Value = '0.1412'
$ws.Cells.Item(25, 5).Value = '  +1.21%  '
$ws.Cells.Item(26, 4).Value = '8.395'
$ws.Cells.Item(26, 5).Value = '  -0.72%  '
$ws.Cells.Item(27, 4).Value = '17.68'
$ws.Cells.Item(27, 5).Value = '  -0.54%  '
$ws.Cells.Item(28, 4).Value = '1.470'
$ws.Cells.Item(28, 5).Value = '  -1.37%  '
$ws.Cells.Item(29, 4).Value = '0.05717'
$ws.Cells.Item(29, 5).Value = '  -2.92%  '
$ws.Cells.Item(30, 4).Value = '1.248'
$ws.Cells.Item(30, 5).Value = '  -2.36%  '
$ws.Cells.Item(31, 4).Value = '4.137'
$ws.Cells.Item(31, 5).Value = '  +0.23%  '
$ws.Cells.Item(32, 4).Value = '4.035'
$ws.Cells.Item(32, 5).Value = '  -0.76%  '
$ws.Cells.Item(33, 5).Value = '  -2.80%  '
$ws.Cells.Item(34, 4).Value = '1.158'
$ws.Cells.Item(34, 5).Value = '  -1.42%  '
$ws.Cells.Item(35, 4).Value = '0.7173'
$ws.Cells.Item(35, 5).Value = '  -0.39%  '
$ws.Cells.Item(36, 4).Value = '2.591'
$ws.Cells.Item(36, 5).Value = '  -0.32%  '
$ws.Cells.Item(37, 4).Value = '1.253.06'
$ws.Cells.Item(37, 5).Value = '  +1.87%  '
$ws.Cells.Item(38, 4).Value = '0.01813'
$ws.Cells.Item(38, 5).Value = '  +1.86%  '
$ws.Cells.Item(39, 4).Value = '2.790'
$ws.Cells.Item(39, 5).Value = '  -0.61%  '
$ws.Cells.Item(40, 4).Value = '0.9082'
$ws.Cells.Item(40, 5).Value = '  -0.63%  '
$ws.Cells.Item(41, 2).Value = 'FraxShare'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(41, 4).Value = '6.149'
$ws.Cells.Item(41, 5).Value = '  +0.45%  '
$ws.Cells.Item(42, 2).Value = 'PaxDollar'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(42, 4).Value = '0.9994'
$ws.Cells.Item(42, 5).Value = '  -0.03%  '
$ws.Cells.Item(43, 2).Value = 'Quant'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(43, 4).Value = '101.67'
$ws.Cells.Item(43, 5).Value = '  -0.32%  '
$ws.Cells.Item(44, 2).Value = 'Aave'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(44, 4).Value = '66.49'
$ws.Cells.Item(44, 5).Value = '  -1.50%  '
$ws.Cells.Item(45, 4).Value = '0.00000000119'
$ws.Cells.Item(45, 5).Value = '  +1.82%  '
$ws.Cells.Item(46, 2).Value = 'Aptos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(46, 4).Value = '7.066'
$ws.Cells.Item(46, 5).Value = '  -4.82%  '
$ws.Cells.Item(47, 4).Value = '9.141'
$ws.Cells.Item(47, 5).Value = '  -0.11%  '
$ws.Cells.Item(48, 2).Value = 'TheSandbox'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(48, 4).Value = '0.4031'
$ws.Cells.Item(48, 5).Value = '  -0.34%  '
$ws.Cells.Item(49, 2).Value = 'RenderToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(49, 4).Value = '1.704'
$ws.Cells.Item(49, 5).Value = '  +0.23%  '
$ws.Cells.Item(50, 2).Value = 'Algorand'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(50, 4).Value = '0.1128'
$ws.Cells.Item(50, 5).Value = '  +0.28%  '
$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(51, 4).Value = '0.05741'
$ws.Cells.Item(51, 5).Value = '  -0.17%  '
